$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    Occurs on every sheet's Status-ish columns (Overview B/C, zh-cn/de-de C).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item(2)
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item(3)
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2) Handback report: populate "Latest Target File" (F) and
#    "Latest Handback File" (G) hyperlinks/display text, and fill in the
#    real "Latest Handback DateTime" (H) values (replacing the zero-date
#    placeholder) for both language sheets.
# ---------------------------------------------------------------------------

# --- zh-cn sheet ---
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8a5bedac9d982a85f28f1c0b218ec39ba9153f24/e2e/35efe67b-7699-461d-923e-8925f6541628.md",
    "",
    "",
    "35efe67b-7699-461d-923e-8925f6541628.md"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71111cf3de36ac8ae9d8af5264ef3971dc536dc9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/35efe67b-7699-461d-923e-8925f6541628.c91227dc15ef442992aeeaf68e97009782a86854.zh-cn.xlf",
    "",
    "",
    "35efe67b-7699-461d-923e-8925f6541628.c91227dc15ef442992aeeaf68e97009782a86854.zh-cn.xlf"
) | Out-Null

$wsZhCn.Range("H2").Value = "2016-03-20 18:35:46"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8a5bedac9d982a85f28f1c0b218ec39ba9153f24/e2e/86b20423-8bb9-4abc-9c66-0eed75ac83da.md",
    "",
    "",
    "86b20423-8bb9-4abc-9c66-0eed75ac83da.md"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71111cf3de36ac8ae9d8af5264ef3971dc536dc9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/86b20423-8bb9-4abc-9c66-0eed75ac83da.f8813fda73210ddd7e97c76b77c62525b5b8a2c9.zh-cn.xlf",
    "",
    "",
    "86b20423-8bb9-4abc-9c66-0eed75ac83da.f8813fda73210ddd7e97c76b77c62525b5b8a2c9.zh-cn.xlf"
) | Out-Null

$wsZhCn.Range("H3").Value = "2016-03-20 18:35:46"

# --- de-de sheet ---
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8a5bedac9d982a85f28f1c0b218ec39ba9153f24/e2e/35efe67b-7699-461d-923e-8925f6541628.md",
    "",
    "",
    "35efe67b-7699-461d-923e-8925f6541628.md"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71111cf3de36ac8ae9d8af5264ef3971dc536dc9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/35efe67b-7699-461d-923e-8925f6541628.c91227dc15ef442992aeeaf68e97009782a86854.de-de.xlf",
    "",
    "",
    "35efe67b-7699-461d-923e-8925f6541628.c91227dc15ef442992aeeaf68e97009782a86854.de-de.xlf"
) | Out-Null

$wsDeDe.Range("H2").Value = "2016-03-20 18:35:52"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8a5bedac9d982a85f28f1c0b218ec39ba9153f24/e2e/86b20423-8bb9-4abc-9c66-0eed75ac83da.md",
    "",
    "",
    "86b20423-8bb9-4abc-9c66-0eed75ac83da.md"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71111cf3de36ac8ae9d8af5264ef3971dc536dc9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/86b20423-8bb9-4abc-9c66-0eed75ac83da.f8813fda73210ddd7e97c76b77c62525b5b8a2c9.de-de.xlf",
    "",
    "",
    "86b20423-8bb9-4abc-9c66-0eed75ac83da.f8813fda73210ddd7e97c76b77c62525b5b8a2c9.de-de.xlf"
) | Out-Null

$wsDeDe.Range("H3").Value = "2016-03-20 18:35:52"

Write-Output "Handback report generated."
